$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells AB1/AC1: bold, size 11 (theme) group-name headers
$ws.Range("AB1").Value = "MAKRO_gruppe"
$ws.Range("AC1").Value = "Egen_gruppe"
$hdr = $ws.Range("AB1:AC1")
$hdr.Font.Size = 11
$hdr.Font.Bold = $true

# Data rows 2-48: AB = MAKRO_gruppe, AC = Egen_gruppe
$ws.Range("AB2").Value = "Total"
$ws.Range("AC2").Value = "Total"
$ws.Range("AB3").Value = "Varer"
$ws.Range("AC3").Value = "Fødevarer"
$ws.Range("AB4").Value = "Varer"
$ws.Range("AC4").Value = "Fødevarer"
$ws.Range("AB5").Value = "Varer"
$ws.Range("AC5").Value = "Fødevarer"
$ws.Range("AB6").Value = "Varer"
$ws.Range("AC6").Value = "Andre varer"
$ws.Range("AB7").Value = "Varer"
$ws.Range("AC7").Value = "Andre varer"
$ws.Range("AB8").Value = "Varer"
$ws.Range("AC8").Value = "Andre varer"
$ws.Range("AB9").Value = "Varer"
$ws.Range("AC9").Value = "Andre varer"
$ws.Range("AB10").Value = "Bolig"
$ws.Range("AC10").Value = "Bolig"
$ws.Range("AB11").Value = "Bolig"
$ws.Range("AC11").Value = "Bolig"
$ws.Range("AB12").Value = "Bolig"
$ws.Range("AC12").Value = "Bolig"
$ws.Range("AB13").Value = "Bolig"
$ws.Range("AC13").Value = "Bolig"
$ws.Range("AB14").Value = "Energi"
$ws.Range("AC14").Value = "Energi"
$ws.Range("AB15").Value = "Varer"
$ws.Range("AC15").Value = "Andre varer"
$ws.Range("AB16").Value = "Varer"
$ws.Range("AC16").Value = "Andre varer"
$ws.Range("AB17").Value = "Varer"
$ws.Range("AC17").Value = "Andre varer"
$ws.Range("AB18").Value = "Varer"
$ws.Range("AC18").Value = "Andre varer"
$ws.Range("AB19").Value = "Varer"
$ws.Range("AC19").Value = "Andre varer"
$ws.Range("AB20").Value = "Varer"
$ws.Range("AC20").Value = "Andre varer"
$ws.Range("AB21").Value = "Varer"
$ws.Range("AC21").Value = "Andre varer"
$ws.Range("AB22").Value = "Tjenester"
$ws.Range("AC22").Value = "Tjenester"
$ws.Range("AB23").Value = "Tjenester"
$ws.Range("AC23").Value = "Tjenester"
$ws.Range("AB24").Value = "Biler"
$ws.Range("AC24").Value = "Biler"
$ws.Range("AB25").Value = "Biler"
$ws.Range("AC25").Value = "Biler"
$ws.Range("AB26").Value = "Tjenester"
$ws.Range("AC26").Value = "Off transport"
$ws.Range("AB27").Value = "Tjenester"
$ws.Range("AC27").Value = "Tjenester"
$ws.Range("AB28").Value = "Varer"
$ws.Range("AC28").Value = "Andre varer"
$ws.Range("AB29").Value = "Tjenester"
$ws.Range("AC29").Value = "Tjenester"
$ws.Range("AB30").Value = "Varer"
$ws.Range("AC30").Value = "Andre varer"
$ws.Range("AB31").Value = "Tjenester"
$ws.Range("AC31").Value = "Tjenester"
$ws.Range("AB32").Value = "Varer"
$ws.Range("AC32").Value = "Andre varer"
$ws.Range("AB33").Value = "Tjenester"
$ws.Range("AC33").Value = "Tjenester"
$ws.Range("AB34").Value = "Varer"
$ws.Range("AC34").Value = "Andre varer"
$ws.Range("AB35").Value = "Turisme"
$ws.Range("AC35").Value = "Turisme"
$ws.Range("AB36").Value = "Tjenester"
$ws.Range("AC36").Value = "Tjenester"
$ws.Range("AB37").Value = "Tjenester"
$ws.Range("AC37").Value = "Tjenester"
$ws.Range("AB38").Value = "Tjenester"
$ws.Range("AC38").Value = "Tjenester"
$ws.Range("AB39").Value = "Tjenester"
$ws.Range("AC39").Value = "Tjenester"
$ws.Range("AB40").Value = "Tjenester"
$ws.Range("AC40").Value = "Tjenester"
$ws.Range("AB41").Value = "Turisme"
$ws.Range("AC41").Value = "Turisme"
$ws.Range("AB42").Value = "Varer"
$ws.Range("AC42").Value = "Andre varer"
$ws.Range("AB43").Value = "Tjenester"
$ws.Range("AC43").Value = "Tjenester"
$ws.Range("AB44").Value = "Varer"
$ws.Range("AC44").Value = "Andre varer"
$ws.Range("AB45").Value = "Tjenester"
$ws.Range("AC45").Value = "Tjenester"
$ws.Range("AB46").Value = "Tjenester"
$ws.Range("AC46").Value = "Tjenester"
$ws.Range("AB47").Value = "Tjenester"
$ws.Range("AC47").Value = "Tjenester"
$ws.Range("AB48").Value = "Tjenester"
$ws.Range("AC48").Value = "Tjenester"

# Update selection to match the authored view (AB1:AC1 active)
[void]$ws.Range("AB1:AC1").Select()
